# Add a comment anchored on the final "e" of the second occurrence of
# "line" in the paragraph "The double cannon is nothing different ...".
# That occurrence is the one in "A <ms>line</ms> is a <ms>twelfth part</ms>
# of a cannonball.", which is being glossed by the new comment.

# New comments created via Comments.Add() pick up their author from the
# current Word "user name" at the moment of creation (just like real Word
# COM automation), so set it before adding the comment.
$word.UserName = "Colin Debuiche"

$d = $word.ActiveDocument

# Locate the unique, unambiguous run of text that contains the specific
# "line" occurrence we need (the one immediately before " is a <ms>twelfth
# part</ms>"), using the literal markup characters that are present in the
# document's own text.
$search = $d.Content
$found = $search.Find.Execute("<ms>line</ms> is a <ms>twelfth part</ms>", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not locate the target 'line' occurrence"
}

# "<ms>" is 4 characters, "line" is the next 4 characters.
$lineStart = $search.Start + 4
$lineEnd = $lineStart + 4
$lineRange = $d.Range($lineStart, $lineEnd)
if ($lineRange.Text -ne "line") {
    throw "Unexpected text at computed range: '$($lineRange.Text)'"
}

# Only the final letter "e" is the commented span.
$eRange = $d.Range($lineEnd - 1, $lineEnd)
if ($eRange.Text -ne "e") {
    throw "Unexpected text at computed 'e' range: '$($eRange.Text)'"
}

$comment = $d.Comments.Add($eRange, "twelfth of a king foot in this context (PB)")
